$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "E005EE04"
$ws.Range("D11").Select()
